$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "45x195 mm. ubh. Spaertrae"
$ws.Range("D7").Value = "19x100 mm. trykimp. Braet"
$ws.Range("D33").Value = "1x20 mm. hulband 10 mtr."
$ws.Range("D20").Value = "universal hojre"
$ws.Range("D28").Value = "10x120 mm. Braeddebolt"
$ws.Range("D22").Value = "50x75 mm. Stalddorsgreb"
$ws.Range("D23").Value = "t-haengsel"
$ws.Range("D2").Value = "25x150 mm. trykimp. Braet"
$ws.Range("D3").Value = "faedigskaret"
$ws.Range("D8").Value = "25x50 mm. trykimp. Braet"
$ws.Range("D17").Value = "B&C Toplaegte holder"
$ws.Range("D11").Value = "25x125 mm. trykimp. Braet"
$ws.Range("D10").Value = "25x200 mm. trykimp. Braet"
$ws.Range("D13").Value = "38x73 mm. Laegte"
$ws.Range("D9").Value = "38x73 mm. Taglaegte T1"

$ws.Range("D20").Select()
